$d = $word.ActiveDocument

# --- 1. Remove the old "_GoBack" bookmark that sits near the top heading
#        (just after the "2008-08" run, before " Commands Spec,") ---
if ($d.Bookmarks.Exists("_GoBack")) {
    [void]$d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Normalize the two date smartTags in the "Date: ..." paragraph so
#        their <w:attr> children are written Year/Day/Month instead of
#        Month/Day/Year (values unchanged: 21 June 2008 - 31 Aug 2008) ---
$dateXml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00CF03FF" w:rsidRPr="00E776B2" w:rsidRDefault="00CF03FF" w:rsidP="00FF5F4B"><w:pPr><w:ind w:left="284"/><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr></w:pPr><w:r w:rsidRPr="00E776B2"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t xml:space="preserve">Date: </w:t></w:r><w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date"><w:smartTagPr><w:attr w:name="Year" w:val="2008"/><w:attr w:name="Day" w:val="21"/><w:attr w:name="Month" w:val="6"/></w:smartTagPr><w:r w:rsidR="00C26530" w:rsidRPr="00E776B2"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t>June 21, 2008</w:t></w:r></w:smartTag><w:r w:rsidR="00C26530" w:rsidRPr="00E776B2"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t xml:space="preserve"> – </w:t></w:r><w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date"><w:smartTagPr><w:attr w:name="Year" w:val="2008"/><w:attr w:name="Day" w:val="31"/><w:attr w:name="Month" w:val="8"/></w:smartTagPr><w:r w:rsidR="00744698"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t xml:space="preserve">August </w:t></w:r><w:r w:rsidR="00F90E05"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t>31</w:t></w:r><w:r w:rsidR="00C26530" w:rsidRPr="00E776B2"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t>, 2008</w:t></w:r></w:smartTag></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$d.Paragraphs.Item(5).Range.InsertXML($dateXml1)

# --- 3. Remove the two empty, indented (ind left=284, italic) paragraphs
#        that sit right before the "Goal" heading ---
[void]$d.Paragraphs.Item(6).Range.Delete()
[void]$d.Paragraphs.Item(6).Range.Delete()

# --- 4. Rewrite the "Goal" heading paragraph: drop the <w:smartTag
#        w:element="place"> wrapper around "Goa" (keep the plain run) and
#        add the "_GoBack" bookmark at the start of the paragraph ---
$goalXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00EF400A" w:rsidRDefault="00EF400A" w:rsidP="00FF5F4B"><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>Goa</w:t></w:r><w:r w:rsidR="00FF5F4B"><w:t>l</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$d.Paragraphs.Item(6).Range.InsertXML($goalXml)

# --- 5. Normalize the two date smartTags in the "Time" paragraph
#        (June 28, 2008 - August 31, 2008) the same way as step 2 ---
$timeXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00FF5F4B" w:rsidRDefault="006460FF" w:rsidP="001A67E0"><w:pPr><w:ind w:left="426"/></w:pPr><w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date"><w:smartTagPr><w:attr w:name="Year" w:val="2008"/><w:attr w:name="Day" w:val="28"/><w:attr w:name="Month" w:val="6"/></w:smartTagPr><w:r><w:t>June 28</w:t></w:r><w:r w:rsidR="00916C2E"><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="00FF5F4B"><w:t>2008</w:t></w:r></w:smartTag><w:r w:rsidR="00916C2E"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00FF5F4B"><w:t xml:space="preserve">– </w:t></w:r><w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date"><w:smartTagPr><w:attr w:name="Year" w:val="2008"/><w:attr w:name="Day" w:val="31"/><w:attr w:name="Month" w:val="8"/></w:smartTagPr><w:r w:rsidR="0007449F"><w:t>August 31</w:t></w:r><w:r w:rsidR="00916C2E"><w:t>, 2008</w:t></w:r></w:smartTag></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$d.Paragraphs.Item(15).Range.InsertXML($timeXml)

Write-Output "edit complete"
